$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.83509999999999
$ws.Range("A6").Value = -22.828
$ws.Range("A7").Value = -21.96180000000002
$ws.Range("B7").Value = 4.647400000000002
$ws.Range("A8").Value = -22.31510000000002
$ws.Range("B11").Value = 5.356699999999997
$ws.Range("B12").Value = 4.5661
$ws.Range("E12").Value = 17.9577
$ws.Range("E13").Value = 16.76670000000001
$ws.Range("E14").Value = 16.89490000000001
$ws.Range("B15").Value = 5.260499999999997
$ws.Range("A16").Value = -21.96100000000001
$ws.Range("E16").Value = 16.72909999999999
$ws.Range("E19").Value = 16.338
$ws.Range("A20").Value = -23.1724
$ws.Range("B20").Value = 5.127599999999997
$ws.Range("E20").Value = 16.5236
$ws.Range("A21").Value = -22.3231
$ws.Range("B21").Value = 5.295699999999998
$ws.Range("B22").Value = 10.1268
$ws.Range("E22").Value = 16.75970000000001
$ws.Range("B23").Value = 9.050800000000002
$ws.Range("A28").Value = -22.19989999999999
$ws.Range("A29").Value = -21.68610000000001
$ws.Range("B29").Value = 5.248000000000001
$ws.Range("A30").Value = -21.81170000000002
$ws.Range("A32").Value = -21.3866
$ws.Range("B34").Value = 9.535200000000009
$ws.Range("E36").Value = 15.7712
$ws.Range("A40").Value = -19.3595
$ws.Range("B42").Value = 10.2049
$ws.Range("B43").Value = 5.603099999999999
$ws.Range("E43").Value = 17.25030000000002
$ws.Range("B44").Value = 5.025000000000005
$ws.Range("B45").Value = 4.850900000000007
$ws.Range("A46").Value = -22.08680000000002
$ws.Range("B46").Value = 4.827800000000005
$ws.Range("E46").Value = 16.83749999999998
$ws.Range("B50").Value = 4.538399999999998
$ws.Range("E50").Value = 16.6055
$ws.Range("A51").Value = -22.22939999999999
$ws.Range("B51").Value = 5.292999999999996
$ws.Range("A52").Value = -22.06029999999999
$ws.Range("A57").Value = -22.80070000000001
$ws.Range("B57").Value = 5.051299999999995
$ws.Range("A59").Value = -22.3229
$ws.Range("A62").Value = -22.27980000000001
$ws.Range("B65").Value = 5.271399999999999
$ws.Range("A66").Value = -21.4836
$ws.Range("B66").Value = 4.778099999999997
$ws.Range("B67").Value = 5.409199999999998
$ws.Range("A73").Value = -19.37339999999999
$ws.Range("A74").Value = -21.88729999999998
$ws.Range("E76").Value = 16.63869999999999
$ws.Range("A77").Value = -20.17219999999998
$ws.Range("B79").Value = 9.599900000000009
$ws.Range("B84").Value = 5.470600000000001
$ws.Range("B87").Value = 5.106399999999997
$ws.Range("A92").Value = -21.35790000000002
$ws.Range("B92").Value = 5.220999999999993
$ws.Range("E95").Value = 18.03110000000001
$ws.Range("B97").Value = 6.044399999999999
$ws.Range("E97").Value = 16.72
$ws.Range("E99").Value = 16.66939999999999
$ws.Range("A100").Value = -22.21829999999999
